$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Preserve the bold/centered/bordered header style that currently
#        lives on A1 (style also shared by A2:A4) so it can be re-applied to
#        the reshaped header row after the rebuild. ---
$ws.Range("A1").Copy()

# --- 2. Wipe every existing cell (content + formatting) so the sheet starts
#        clean before we rebuild it with the new shape. ---
$ws.UsedRange.Clear()

# --- 3. The "year" / numeric-looking values must be written as literal TEXT
#        (matching the source data's shared-string storage), not numbers, so
#        mark the destination cells as Text before entering them. ---
$ws.Range("B2:D2").NumberFormat = "@"
$ws.Range("B4:D6").NumberFormat = "@"
$ws.Range("B8:D9").NumberFormat = "@"

# ---- Row 1: new column headers ----
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"

# ---- Row 2: units / year header ----
$ws.Range("A2").Value = "Kilowatt hours (kWh)"
$ws.Range("B2").Value = "2019"
$ws.Range("C2").Value = "2018"
$ws.Range("D2").Value = "2017"

# ---- Electricity consumption section ----
$ws.Range("A3").Value = "Electricity consumption"

$ws.Range("A4").Value = "Total electricity consumption"
$ws.Range("B4").Value = "33,104,461"
$ws.Range("C4").Value = "33,035,150"
$ws.Range("D4").Value = "32,208,132"

$ws.Range("A5").Value = "of which green electricity"
$ws.Range("B5").Value = "32,782,553"
$ws.Range("C5").Value = "33,005,705"
$ws.Range("D5").Value = "32,036,926"

$ws.Range("A6").Value = "Individual electricity consumption (kWh per FTE)"
$ws.Range("B6").Value = "7,116"
$ws.Range("C6").Value = "7,000"
$ws.Range("D6").Value = "6,734"

# ---- District heating consumption section ----
$ws.Range("A7").Value = "District heating consumption"

$ws.Range("A8").Value = "Total district heating consumption"
$ws.Range("B8").Value = "18,964,126"
$ws.Range("C8").Value = "18,124,104"
$ws.Range("D8").Value = "18,563,309"

$ws.Range("A9").Value = "Individual heating consumption (kWh per FTE)"
$ws.Range("B9").Value = "4,077"
$ws.Range("C9").Value = "3,841"
$ws.Range("D9").Value = "3,881"

# --- 4. Now that the text is safely entered, drop the temporary Text
#        number-format again so the data cells are back to the plain/
#        default look (no explicit format) like the rest of the sheet. ---
$ws.Range("B2:D2").ClearFormats()
$ws.Range("B4:D6").ClearFormats()
$ws.Range("B8:D9").ClearFormats()

# --- 5. Re-apply the preserved header style to the full new header row
#        (including A1, which got wiped by UsedRange.Clear() above). ---
$ws.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
